$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B to hold the new "Jun_26" week.
# This shifts the existing date columns (old B,C,D,E -> new C,D,E,F).
$ws.Columns("B:B").Insert()

# New header for the freshly inserted week column.
$ws.Cells.Item(1, 2).Value = "Jun_26"

# Fill the new week's column with the "no change" placeholder used
# throughout the rest of the sheet for every existing analyst row.
for ($r = 2; $r -le 27; $r++) {
    $ws.Cells.Item($r, 2).Value = "UN"
}

# Two new analyst rows added at the bottom of the table.
$ws.Cells.Item(28, 1).Value = "Benchmark"
$ws.Cells.Item(28, 2).Value = "UN"

$ws.Cells.Item(29, 1).Value = "Evercore ISI"
$ws.Cells.Item(29, 2).Value = "UN"
